$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the duplicate bold "Play Divina Commedia Free..." paragraph
#    that sits right before the final (italic) paragraph, and replace
#    the final paragraph's text with the new "Create a Cartoon..." copy
#    (keeping its italic run formatting intact).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldDupPara = $d.Paragraphs.Item($count - 1)
$finalPara = $d.Paragraphs.Item($count)

# Delete the whole bold-duplicate paragraph, including its paragraph mark.
$delRange = $d.Range($boldDupPara.Range.Start, $finalPara.Range.Start)
$delRange.Delete()

# Replace the (now last) italic paragraph's text in place.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$italicRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$italicRange.Text = "Create a Cartoon Style Feature Image for Divina Commedia Featuring a Happy Maya Warrior with Glasses: For the feature image of Divina Commedia, we want a fun and cartoony style to match the game's graphics. To add an exciting element, we want to showcase a Maya warrior with glasses who's thrilled to be exploring the depths of Hell. The image could show the Maya warrior grinning with delight as he spins the reels and wins rewards. He could be in a typical Mayan warrior outfit, complete with a headdress and chest armor. His glasses could be comically large and could reflect the fire and brimstone of Hell behind him. We could also include some of the game's symbols in the image, such as the gold coin, pitchfork, and Roman numerals, to give viewers a glimpse of what they can expect from the game. Overall, we want the feature image to be fun, lighthearted, and engaging, showcasing the excitement and adventure of playing Divina Commedia."

# ---------------------------------------------------------------------
# 2) Insert the new "Meta description" paragraph right after the title
#    (first) paragraph.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# Type the non-bold suffix first so the insertion point for the bold
# lead-in stays at the very start of the paragraph.
$suffixRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start)
$suffixRange.Text = ": Read our review of Divina Commedia, a fun and humorous video slot inspired by Dante's famous Italian work. Play free and discover two exciting bonus games and four jackpots."

$leadRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start)
$leadRange.Text = "Meta description"

$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + 16)
$boldRange.Bold = $true
